# Recipe_Scrapping_Data.xlsx edit:
#  - Duplicate the label column (column A) of the "LFV" sheet into a new
#    sheet named "LCHF", placed right after "LFV".
#  - "LCHF" keeps only the label column (no ingredient/recipe detail
#    columns), gets its own wider column A, and becomes the active sheet.
#  - The original "LFV" sheet is no longer the active tab/selection.

$wb  = $excel.ActiveWorkbook
$lfv = $wb.Worksheets.Item(1)

# Add a new, blank worksheet right after "LFV" and give it its final name.
$lchf = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lfv)
$lchf.Name = "LCHF"

# Bring over the label column (values + styles) from LFV's column A.
$lfv.Range("A1:A14").Copy($lchf.Range("A1"))

# Row heights for the two taller wrapped-text rows (not carried by Copy).
$lchf.Rows("3:3").RowHeight = 39.4
$lchf.Rows("4:4").RowHeight = 26.25

# Widen column A on the new sheet.
$lchf.Columns("A:A").ColumnWidth = 34.3

# Update selections: LFV no longer holds the active selection/tab;
# LCHF becomes the active sheet with C8 selected.
$lfv.Range("A1:A14").Select()
$lchf.Range("C8").Select()
$lchf.Activate()
